$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max" header / numeric value in row2) entirely,
# shifting columns D and E left by one (D->C, E->D).
$ws.Range("C1:C2").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Update the remaining B2 value to the new prediction number.
$ws.Range("B2").Value = 45680239.69977608
